# Auto-generated edits applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.987.92"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "3.397.61"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'573.26"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "'142.78"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("D7").Value = "3.399.09"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "3.975.77"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D15").Value = "'28.06"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "3.393.61"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "61.048.49"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D20").Value = "'13.86"
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("E21").Value = "  -4.68%  "
$ws.Range("D22").Value = "'383.91"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "'74.24"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -4.64%  "
$ws.Range("D27").Value = "3.533.82"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'23.53"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "'7.00"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "'167.73"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "3.427.70"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0773"
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'27.61"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").Value = "'0.782"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").Value = "'1.14"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").Value = "2.485.80"
$ws.Range("E48").Value = "  -4.86%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").Value = "'0.0266"
$ws.Range("E51").Value = "  +0.92%  "
